$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1.73
$ws.Range("M2").Value = 2
$ws.Range("P2").Value = 1.78
$ws.Range("Q2").Value = 2.03
$ws.Range("K3").Value = 7.5
$ws.Range("AD3").Value = 501
$ws.Range("AE3").Value = 8
$ws.Range("G4").Value = 2.6
$ws.Range("I4").Value = 2.7
$ws.Range("AB4").Value = 15
$ws.Range("AI4").Value = 23
$ws.Range("G6").Value = 4.5
$ws.Range("I6").Value = 1.73
$ws.Range("L6").Value = 1.2
$ws.Range("M6").Value = 4.33
$ws.Range("R6").Value = 1.67
$ws.Range("S6").Value = 2.1
$ws.Range("U6").Value = 26
$ws.Range("W6").Value = 51
$ws.Range("X6").Value = 34
$ws.Range("AF6").Value = 9
$ws.Range("K7").Value = 12
$ws.Range("Z7").Value = 12
$ws.Range("AE7").Value = 11
$ws.Range("P8").Value = 1.75
$ws.Range("Q8").Value = 2.05
$ws.Range("AA8").Value = 6
$ws.Range("AB8").Value = 23
$ws.Range("AG8").Value = 17
$ws.Range("J9").Value = 1.13
$ws.Range("K9").Value = 6
$ws.Range("P9").Value = 1.57
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 2.1
$ws.Range("S9").Value = 1.67
$ws.Range("V9").Value = 10
$ws.Range("AE10").Value = 10
$ws.Range("AG10").Value = 17
$ws.Range("G11").Value = 1.9
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 3.9
$ws.Range("N11").Value = 2.2
$ws.Range("O11").Value = 1.65
$ws.Range("R11").Value = 1.91
$ws.Range("S11").Value = 1.8
$ws.Range("AE11").Value = 10
$ws.Range("AG11").Value = 15
$ws.Range("J12").Value = 1.05
$ws.Range("K12").Value = 11
$ws.Range("Y12").Value = 29
$ws.Range("AD12").Value = 351
$ws.Range("AG12").Value = 21
$ws.Range("J13").Value = 1.08
$ws.Range("K13").Value = 8
$ws.Range("T13").Value = 9
$ws.Range("Z13").Value = 8
$ws.Range("G15").Value = 3.1
$ws.Range("H15").Value = 3.4
$ws.Range("L15").Value = 1.29
$ws.Range("M15").Value = 3.5
$ws.Range("N15").Value = 1.98
$ws.Range("O15").Value = 1.88
$ws.Range("R15").Value = 1.73
$ws.Range("S15").Value = 2
$ws.Range("U15").Value = 17
$ws.Range("V15").Value = 11
$ws.Range("AB15").Value = 13
$ws.Range("AC15").Value = 41
$ws.Range("AD15").Value = 201
$ws.Range("AI15").Value = 17
$ws.Range("AJ15").Value = 26
$ws.Range("H16").Value = 6.6
$ws.Range("I16").Value = 1.09
$ws.Range("N16").Value = 1.44
$ws.Range("O16").Value = 2.6
$ws.Range("R16").Value = 2.7
$ws.Range("S16").Value = 1.4
$ws.Range("T16").Value = 50
$ws.Range("Z16").Value = 14.5
$ws.Range("AA16").Value = 14.5
$ws.Range("AB16").Value = 40
$ws.Range("AE16").Value = 6.6
$ws.Range("AF16").Value = 4.75
$ws.Range("AG16").Value = 10
$ws.Range("AH16").Value = 4.8
$ws.Range("AI16").Value = 10
$ws.Range("AJ16").Value = 37
$ws.Range("H17").Value = 3.55
$ws.Range("L17").Value = 1.2
$ws.Range("M17").Value = 4
$ws.Range("R17").Value = 1.59
$ws.Range("S17").Value = 2.22
$ws.Range("U17").Value = 8.25
$ws.Range("AA17").Value = 6.2
$ws.Range("AF17").Value = 16.5
$ws.Range("AH17").Value = 37
$ws.Range("G18").Value = 2.2
$ws.Range("I18").Value = 3.2
$ws.Range("K18").Value = 12
$ws.Range("T18").Value = 9
$ws.Range("V18").Value = 9
$ws.Range("AG18").Value = 12
$ws.Range("G19").Value = 1.53
$ws.Range("H19").Value = 4.2
$ws.Range("I19").Value = 6.25
$ws.Range("J19").Value = 1.07
$ws.Range("K19").Value = 9
$ws.Range("R19").Value = 2.2
$ws.Range("S19").Value = 1.62
$ws.Range("W19").Value = 10
$ws.Range("AB19").Value = 23
$ws.Range("AE19").Value = 13
$ws.Range("L20").Value = 1.36
$ws.Range("M20").Value = 3
$ws.Range("N20").Value = 2.2
$ws.Range("O20").Value = 1.65
$ws.Range("R20").Value = 1.91
$ws.Range("S20").Value = 1.8
$ws.Range("T20").Value = 7.5
$ws.Range("AB20").Value = 17
$ws.Range("R21").Value = 2.1
$ws.Range("S21").Value = 1.67
$ws.Range("X21").Value = 34
$ws.Range("AF21").Value = 9.5
